# Edit workbook per commit: "New files, changes to files and pipeline after
# bringing in memory optimized code."

$wb = $excel.ActiveWorkbook

# --- Sheet "Trends Status": Insufficient Data row count 184 -> 185 ---
$trendsStatus = $wb.Worksheets.Item("Trends Status")
$trendsStatus.Range("B8").Value = 185
$trendsStatus.Range("C8").Value = 185

# --- Sheet "Priority Status": High/Moderate/Low counts updated ---
$priorityStatus = $wb.Worksheets.Item("Priority Status")
$priorityStatus.Range("B2").Value = 103
$priorityStatus.Range("B3").Value = 286
$priorityStatus.Range("B4").Value = 554

# --- Sheet "Species qualification": label + selected-for-analysis count ---
$speciesQual = $wb.Worksheets.Item("Species qualification")
$speciesQual.Range("A2").Value = "SoIB Assessment"
$speciesQual.Range("B2").Value = 185

# --- Capture the existing "High Priority break-up" sheet's data before we
#     overwrite it, so it can be copied to the brand-new sheet below. ---
$oldBreakUp = $wb.Worksheets.Item("High Priority break-up")

# Rename the existing sheet and replace its contents with the new
# "Interannual update" figures.
$oldBreakUp.Name = "Interannual update - High Pri"

$oldBreakUp.Range("A1").Value = "Break-up"
$oldBreakUp.Range("B1").Value = "High Species (no.)"
$oldBreakUp.Range("C1").Value = "High Species (perc.)"
$oldBreakUp.Range("D1").Value = "New High Species (no.)"
$oldBreakUp.Range("E1").Value = "New High Species (perc.)"

$oldBreakUp.Range("A2").Value = "Trend New"
$oldBreakUp.Range("B2").Value = 97
$oldBreakUp.Range("C2").Value = 94.2
$oldBreakUp.Range("D2").Value = 97
$oldBreakUp.Range("E2").Value = 95.09999999999999

$oldBreakUp.Range("A3").Value = "IUCN"
$oldBreakUp.Range("B3").Value = 6
$oldBreakUp.Range("C3").Value = 5.8
$oldBreakUp.Range("D3").Value = 5
$oldBreakUp.Range("E3").Value = 4.9

# --- Add a brand-new sheet at the end holding the original
#     "High Priority break-up" figures, renamed "Major update - High Priority ". ---
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Major update - High Priority "

# Reuse the bold/centered header formatting already used on the other sheets
# (copy format only) instead of re-deriving it, so the style table stays
# identical to the rest of the workbook.
$trendsStatus.Range("A1:E1").Copy()
$newSheet.Range("A1:E1").PasteSpecial(-4122)

$newSheet.Range("A1").Value = "Break-up"
$newSheet.Range("B1").Value = "High Species (no.)"
$newSheet.Range("C1").Value = "High Species (perc.)"
$newSheet.Range("D1").Value = "New High Species (no.)"
$newSheet.Range("E1").Value = "New High Species (perc.)"

$newSheet.Range("A2").Value = "IUCN"
$newSheet.Range("B2").Value = 1
$newSheet.Range("C2").Value = 100
$newSheet.Range("D2").Value = 1
$newSheet.Range("E2").Value = 100

# Restore the originally-selected tab (first sheet) so tabSelected / activeTab
# stay where they were before our edits (Add() above shifts focus).
$trendsStatus.Activate()
